# Generate Report for Handoff
#
# This mirrors a "re-handoff" of 4 files (rows 4-7 in the zh-cn / de-de
# localization tables) whose Priority moved from "low" to "ht" (high),
# and (for zh-cn only) got a fresh "Latest Handoff Datetime" /
# "Latest HO Xliff Generate Date" timestamp.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- zh-cn sheet: rows 4-7 ---
# Column E = Priority: "low" -> "ht"
# Column H = Latest Handoff Datetime: "2016-10-13 14:35:24" -> "2016-10-13 14:36:00"
foreach ($r in 4..7) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-10-13 14:36:00"
}

# --- de-de sheet: rows 4-7 ---
# Column E = Priority: "low" -> "ht"
foreach ($r in 4..7) {
    $wsDeDe.Range("E$r").Value = "ht"
}

# --- Overview sheet: rows 4-7 ---
# Column G = Latest HO Xliff Generate Date: "2016-10-13 14:35:35" -> "2016-10-13 14:36:12"
foreach ($r in 4..7) {
    $wsOverview.Range("G$r").Value = "2016-10-13 14:36:12"
}
